$wb = $excel.ActiveWorkbook

# 1. Insert a new "item_num" column (O) on every sheet that shares the
#    common "Measures"-style header layout, pushing the old "comment"
#    column from O to P.
$sheetsWithHeader = @("Measures","ID","Dems","Dates","NewVars")
foreach ($name in $sheetsWithHeader) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(15).Insert()
    $ws.Range("O1").Value = "item_num"
}

# 2. Fill the new item_num column with 1 for every data row on the
#    Measures sheet (rows 2-81).
$measures = $wb.Worksheets.Item("Measures")
for ($r = 2; $r -le 81; $r++) {
    $measures.Cells.Item($r, 15).Value = 1
}

# 3. Update the worksheet-scoped _FilterDatabase defined name so it still
#    spans the full header row after the new column was inserted.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -like "*_FilterDatabase*") {
        $nm.RefersTo = "=Measures!`$A`$1:`$P`$1"
    }
}

# 4. Turn off the AutoFilter dropdowns on the Measures sheet.
$measures.AutoFilterMode = $false

# 5. Restore each sheet's selection to mirror the authored workbook state.
$idSheet = $wb.Worksheets.Item("ID")
$idSheet.Range("O1:O1048576").Select()

$demsSheet = $wb.Worksheets.Item("Dems")
$demsSheet.Range("O1:O1048576").Select()

$datesSheet = $wb.Worksheets.Item("Dates")
$datesSheet.Range("O1:O1048576").Select()

$newVarsSheet = $wb.Worksheets.Item("NewVars")
$newVarsSheet.Range("O1:O1048576").Select()

# 6. Make Measures the active sheet/tab, with the same frozen-pane
#    selection the author ended up with.
$measures.Activate()
$measures.Range("O5:O81").Select()

Write-Host "edit applied"
